$wb = $excel.ActiveWorkbook

# Sheet ALC, row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 63977.715
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 74557.336
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 223672.008
$ws.Range("M6").Value = -1388
$ws.Range("N6").Value = -223896.008

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2210
$ws.Range("I98").Value = 2212.5
$ws.Range("J98").Value = 2200
$ws.Range("K98").Value = 2212.5
$ws.Range("L98").Value = 2200
$ws.Range("M98").Value = -714.5
$ws.Range("N98").Value = -5196

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3650
$ws.Range("I113").Value = 3850
$ws.Range("K113").Value = 3850
$ws.Range("M113").Value = -596

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2210
$ws.Range("I122").Value = 2212.5
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 6637.5
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -4187.5
$ws.Range("N122").Value = -11500

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 885.29266
$ws.Range("J129").Value = 955.8857400000001
$ws.Range("L129").Value = 2867.65722
$ws.Range("N129").Value = -12867.65722

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3574180.2
$ws.Range("I137").Value = 7694877
$ws.Range("J137").Value = 2910.0667
$ws.Range("K137").Value = 23084631
$ws.Range("L137").Value = 8730.2001
$ws.Range("M137").Value = -23082081
$ws.Range("N137").Value = -13830.2001

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2152891.5
$ws.Range("I138").Value = 73087.28999999999
$ws.Range("J138").Value = 2607848.8
$ws.Range("K138").Value = 219261.87
$ws.Range("L138").Value = 7823546.399999999
$ws.Range("M138").Value = -214121.87
$ws.Range("N138").Value = -7833826.399999999

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4288583.5
$ws.Range("I32").Value = 4847366
$ws.Range("J32").Value = 4583.5557
$ws.Range("K32").Value = 4847366
$ws.Range("L32").Value = 4583.5557
$ws.Range("M32").Value = -4847079
$ws.Range("N32").Value = -5157.5557

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 77077960
$ws.Range("I61").Value = 100100904
$ws.Range("J61").Value = 334838
$ws.Range("K61").Value = 100100904
$ws.Range("L61").Value = 334838
$ws.Range("M61").Value = -100100692
$ws.Range("N61").Value = -335262

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10957671
$ws.Range("I74").Value = 17929212
$ws.Range("J74").Value = 113050
$ws.Range("K74").Value = 17929212
$ws.Range("L74").Value = 113050
$ws.Range("M74").Value = -17928338
$ws.Range("N74").Value = -114798

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10957671
$ws.Range("I77").Value = 17929212
$ws.Range("J77").Value = 113050
$ws.Range("K77").Value = 89646060
$ws.Range("L77").Value = 565250
$ws.Range("M77").Value = -89641692
$ws.Range("N77").Value = -573986

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2000758.6
$ws.Range("I110").Value = 3333833.2
$ws.Range("J110").Value = 1146.5
$ws.Range("K110").Value = 3333833.2
$ws.Range("L110").Value = 1146.5
$ws.Range("M110").Value = -3331788.2
$ws.Range("N110").Value = -5236.5

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 77077960
$ws.Range("I136").Value = 100100904
$ws.Range("J136").Value = 334838
$ws.Range("K136").Value = 300302712
$ws.Range("L136").Value = 1004514
$ws.Range("M136").Value = -300300162
$ws.Range("N136").Value = -1009614

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3305
$ws.Range("I107").Value = 2988
$ws.Range("J107").Value = 3833.3333
$ws.Range("K107").Value = 2988
$ws.Range("L107").Value = 3833.3333
$ws.Range("M107").Value = -1068
$ws.Range("N107").Value = -7673.3333

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 831.22
$ws.Range("I31").Value = 768.5
$ws.Range("J31").Value = 844.9878
$ws.Range("K31").Value = 768.5
$ws.Range("L31").Value = 844.9878
$ws.Range("M31").Value = -473.5
$ws.Range("N31").Value = -1434.9878

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 831.22
$ws.Range("I34").Value = 768.5
$ws.Range("J34").Value = 844.9878
$ws.Range("K34").Value = 768.5
$ws.Range("L34").Value = 844.9878
$ws.Range("M34").Value = -566.5
$ws.Range("N34").Value = -1248.9878

# Sheet CRP, row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7500
$ws.Range("I99").Value = 7285.7144
$ws.Range("J99").Value = 9000
$ws.Range("K99").Value = 7285.7144
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -5787.7144
$ws.Range("N99").Value = -11996

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 535
$ws.Range("I107").Value = 404.6842
$ws.Range("J107").Value = 741.3333
$ws.Range("K107").Value = 404.6842
$ws.Range("L107").Value = 741.3333
$ws.Range("M107").Value = 1515.3158
$ws.Range("N107").Value = -4581.3333

# Sheet CRP, row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7500
$ws.Range("I126").Value = 7285.7144
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 21857.1432
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -19387.1432
$ws.Range("N126").Value = -31940

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 122015.11
$ws.Range("I134").Value = 1908
$ws.Range("J134").Value = 182068.67
$ws.Range("K134").Value = 5724
$ws.Range("L134").Value = 546206.01
$ws.Range("M134").Value = -3189
$ws.Range("N134").Value = -551276.01

# Sheet CUL, row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 785.5909
$ws.Range("I68").Value = 450.14285
$ws.Range("J68").Value = 942.13336
$ws.Range("K68").Value = 1350.42855
$ws.Range("L68").Value = 2826.40008
$ws.Range("M68").Value = -539.4285500000001
$ws.Range("N68").Value = -4448.40008

# Sheet CUL, row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 785.5909
$ws.Range("I71").Value = 450.14285
$ws.Range("J71").Value = 942.13336
$ws.Range("K71").Value = 4051.28565
$ws.Range("L71").Value = 8479.20024
$ws.Range("M71").Value = 4.71434999999974
$ws.Range("N71").Value = -16591.20024

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 953
$ws.Range("I131").Value = 485.55554
$ws.Range("J131").Value = 1303.5834
$ws.Range("K131").Value = 1456.66662
$ws.Range("L131").Value = 3910.7502
$ws.Range("M131").Value = 3583.33338
$ws.Range("N131").Value = -13990.7502

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1732.02
$ws.Range("I140").Value = 797.5
$ws.Range("J140").Value = 1910.0238
$ws.Range("K140").Value = 2392.5
$ws.Range("L140").Value = 5730.0714
$ws.Range("M140").Value = 2787.5
$ws.Range("N140").Value = -16090.0714
